# Update "想去人数" (attendance count) figures on two worksheets that
# both list the same set of events: "展览" (Exhibitions) and "全部类型" (All types).
#
# 展览 sheet:
#   F2: 14331 -> 14381
#   F5: 238   -> 239
#   F6: 567   -> 571
#   F7: 1504  -> 1512
#
# 全部类型 sheet (same events, but at different rows because this sheet also
# interleaves 演出/本地生活 entries):
#   F2: 14331 -> 14381
#   F5: 238   -> 239
#   F8: 567   -> 571
#   F9: 1504  -> 1512

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 14381
$wsExhibitions.Range("F5").Value = 239
$wsExhibitions.Range("F6").Value = 571
$wsExhibitions.Range("F7").Value = 1512

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 14381
$wsAllTypes.Range("F5").Value = 239
$wsAllTypes.Range("F8").Value = 571
$wsAllTypes.Range("F9").Value = 1512
